$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values for rows 2-10 (existing rows) plus new row 11
$data = @(
    @(6,  0.4985553813290332, 2.131733426832194, 17.05044517982819, 4.129218470828128, 4.139797756631062, 51),
    @(7,  0.3670740509330708, 2.310008463032842, 18.22424447820178, 4.268986352543397, 4.296355968907215, 50),
    @(8,  0.6013248752497656, 2.260079145283003, 17.85224903368007, 4.225192188963725, 4.225523177608233, 49),
    @(9,  0.3757086490286848, 2.257680262786419, 19.07933980741992, 4.36799036256033,  4.397854349388155, 48),
    @(10, 0.6368657846289102, 2.180689379475908, 18.26526036982974, 4.27378759063079,  4.271757936889329, 47),
    @(11, 0.2352128511284962, 1.884428425708629, 15.78161112617242, 3.972607597809331, 4.009458699620951, 46),
    @(12, 0.4131478667422561, 1.866528751697462, 15.52343479139494, 3.939979034385201, 3.9625332330756,   45),
    @(13, 0.2445497104634806, 1.984813716461438, 16.20214174503358, 4.02518841112234,  4.064202304816381, 44),
    @(14, 0.4887008263018376, 1.804638728127334, 15.84511544644097, 3.980592348688945, 3.997232193593323, 43),
    @(15, 0.2035882577814556, 1.941069045450287, 16.47627633247782, 4.059097970298058, 4.103130222273781, 42)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
    $ws.Cells.Item($row, 6).Value = $entry[5]
    $ws.Cells.Item($row, 7).Value = $entry[6]
    $row++
}

# New row 11, column A = "Q9" with same style as other A cells (copy from A10)
$ws.Range("A11").Value = "Q9"
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122) # xlPasteFormats
